# Update "想去人数" (number of people wanting to go) counts on the
# "展览" and "全部类型" sheets:
#   F3: 83 -> 85
#   F4: 52 -> 54

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 85
    $ws.Range("F4").Value = 54
}
